$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns contain values that look numeric (e.g. "1.00",
# "0.0000190", "603.81") but must be preserved as literal text, matching
# the original workbook's inline-string cells. Pre-formatting the range as
# text prevents Excel's COM layer from silently coercing these strings into
# numbers (which would drop formatting like trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value2 = '65.206.43'
$ws.Range("E2").Value2 = '  +2.56%  '
$ws.Range("D3").Value2 = '2.643.63'
$ws.Range("E3").Value2 = '  +1.43%  '
$ws.Range("E4").Value2 = '  +0.07%  '
$ws.Range("D5").Value2 = '601.87'
$ws.Range("E5").Value2 = '  +1.91%  '
$ws.Range("D6").Value2 = '156.53'
$ws.Range("E6").Value2 = '  +4.54%  '
$ws.Range("E7").Value2 = '  +0.03%  '
$ws.Range("E8").Value2 = '  +0.80%  '
$ws.Range("E9").Value2 = '  +10.88%  '
$ws.Range("D10").Value2 = '0.411'
$ws.Range("E10").Value2 = '  +6.13%  '
$ws.Range("D11").Value2 = '5.81'
$ws.Range("E11").Value2 = '  +0.96%  '
$ws.Range("E12").Value2 = '  +2.30%  '
$ws.Range("D13").Value2 = '29.39'
$ws.Range("E13").Value2 = '  +6.27%  '
$ws.Range("D14").Value2 = '0.0000190'
$ws.Range("E14").Value2 = '  +21.61%  '
$ws.Range("D15").Value2 = '3.120.12'
$ws.Range("E15").Value2 = '  +1.56%  '
$ws.Range("D16").Value2 = '65.049.34'
$ws.Range("E16").Value2 = '  +2.66%  '
$ws.Range("D17").Value2 = '2.647.51'
$ws.Range("E17").Value2 = '  +2.00%  '
$ws.Range("D18").Value2 = '12.67'
$ws.Range("E18").Value2 = '  +4.64%  '
$ws.Range("E19").Value2 = '  +3.85%  '
$ws.Range("D20").Value2 = '359.73'
$ws.Range("E20").Value2 = '  +4.12%  '
$ws.Range("E21").Value2 = '  +7.57%  '
$ws.Range("D22").Value2 = '1.00'
$ws.Range("E22").Value2 = '  -0.01%  '
$ws.Range("D23").Value2 = '69.34'
$ws.Range("E23").Value2 = '  +3.78%  '
$ws.Range("D24").Value2 = '1.71'
$ws.Range("E24").Value2 = '  +0.94%  '
$ws.Range("D25").Value2 = '9.44'
$ws.Range("E25").Value2 = '  +2.07%  '
$ws.Range("E26").Value2 = '  +0.06%  '
$ws.Range("D27").Value2 = '8.32'
$ws.Range("E27").Value2 = '  +0.00%  '
$ws.Range("E28").Value2 = '  +2.92%  '
$ws.Range("D29").Value2 = '0.0₃0971'
$ws.Range("E29").Value2 = '  +12.63%  '
$ws.Range("D30").Value2 = '550.97'
$ws.Range("E30").Value2 = '  +0.33%  '
$ws.Range("E31").Value2 = '  +8.85%  '
$ws.Range("E32").Value2 = '  +0.20%  '
$ws.Range("E33").Value2 = '  +1.81%  '
$ws.Range("D34").Value2 = '5.65'
$ws.Range("E34").Value2 = '  +6.10%  '
$ws.Range("D35").Value2 = '6.39'
$ws.Range("E35").Value2 = '  +4.77%  '
$ws.Range("D36").Value2 = '0.431'
$ws.Range("E36").Value2 = '  +4.67%  '
$ws.Range("D37").Value2 = '20.52'
$ws.Range("E37").Value2 = '  +5.53%  '
$ws.Range("E38").Value2 = '  +3.38%  '
$ws.Range("D39").Value2 = '162.42'
$ws.Range("E39").Value2 = '  -2.21%  '
$ws.Range("E40").Value2 = '  +0.15%  '
$ws.Range("E41").Value2 = '  +0.02%  '
$ws.Range("D42").Value2 = '42.81'
$ws.Range("E42").Value2 = '  +7.88%  '
$ws.Range("D43").Value2 = '167.14'
$ws.Range("E43").Value2 = '  +1.24%  '
$ws.Range("D44").Value2 = '4.21'
$ws.Range("E44").Value2 = '  +4.32%  '
$ws.Range("D45").Value2 = '0.0623'
$ws.Range("E45").Value2 = '  +7.43%  '
$ws.Range("D46").Value2 = '2.31'
$ws.Range("E46").Value2 = '  +9.76%  '
$ws.Range("D47").Value2 = '23.27'
$ws.Range("E47").Value2 = '  +1.33%  '
$ws.Range("B48").Value2 = 'Mantle'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value2 = '0.656'
$ws.Range("E48").Value2 = '  +3.68%  '
$ws.Range("B49").Value2 = 'VeChain'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value2 = '0.0263'
$ws.Range("E49").Value2 = '  +5.11%  '
$ws.Range("D50").Value2 = '0.0982'
$ws.Range("E50").Value2 = '  +2.32%  '
$ws.Range("D51").Value2 = '19.68'
$ws.Range("E51").Value2 = '  +2.92%  '
